$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("diagnostics")

$ws.Range("L2").Value = "09b7334a-1514-4de2-b36e-671858821b4f"
$ws.Range("L3").Value = "257bc1e3-35ac-4815-89e5-95244e2d51d6"
$ws.Range("L4").Value = "2310954b-7aa5-44c1-9b20-e3a776857717"
$ws.Range("L5").Value = "91b6a378-6b37-4b75-ae5a-b196a7fd6799"
$ws.Range("L6").Value = "d93cfdb2-6a85-44c9-bc64-02468b7f1b1f"
$ws.Range("L7").Value = "8a286c61-e57a-4c2e-90bf-34b068be4bb6"
$ws.Range("L8").Value = "92d6c614-7196-4c63-80f3-137ae4f34c30"
$ws.Range("L9").Value = "03b6641e-17ae-428e-a19b-ed9fa33300a0"
$ws.Range("L10").Value = "ffe9c8e8-8b86-4473-b860-922d2361d16d"
$ws.Range("L11").Value = "8ba636fc-b681-4c2a-80d6-dd0219d2cd7e"
$ws.Range("L12").Value = "4356e83f-d708-496a-8c6b-ffe83779f9b6"
$ws.Range("L13").Value = "ffa827d9-5f1e-4b91-a168-aeb0b13bb274"
$ws.Range("L14").Value = "c7999171-b179-445d-bc46-bd194d8907e9"
$ws.Range("L15").Value = "11820fe1-94ee-4b33-b0d3-c52c662cf48f"
$ws.Range("L16").Value = "6a22bbc2-3923-46e8-a348-18c416fc2ff6"
$ws.Range("L17").Value = "80fac7d0-28db-4e02-9589-0fa17257fb7c"
$ws.Range("L18").Value = "678fa32b-9d2e-4ef6-84cd-7ce6d01e867f"
$ws.Range("L19").Value = "2a1c34ff-dd37-4ee8-b07e-a22b2cdd7346"
$ws.Range("L20").Value = "4662a0a0-31ee-4be0-8ce3-9b38a410f956"
$ws.Range("L21").Value = "f15353f2-2ee8-41a8-b605-4e8c883cc35e"
$ws.Range("L22").Value = "bea12b49-8b58-48da-ae40-14950f146a97"
$ws.Range("L23").Value = "4570ef72-4d1d-4bbb-ac2a-1a5b8b0bcf43"
$ws.Range("L24").Value = "564795d4-2594-420b-9026-40c350dc8407"
$ws.Range("L25").Value = "8bb86bb1-8d72-44d6-a560-25141a883036"
$ws.Range("L26").Value = "99b95f78-9ad8-42a9-ad5a-44484fb0182a"
$ws.Range("L27").Value = "5778974f-36ac-4681-b599-a4385343977a"
$ws.Range("L28").Value = "c0c70243-00b8-4411-8d91-8d9f38f93bed"
$ws.Range("L29").Value = "1db43e86-ff82-4edc-86fd-bb4ec30ca66e"
$ws.Range("L30").Value = "a1e856ec-ccd9-4831-9ff0-f07042ed67fa"
$ws.Range("L31").Value = "8cdb8af7-12e0-40a7-92f4-d198df058427"
$ws.Range("L32").Value = "9c2eb6c8-384b-432d-a6ef-d8403771d2a2"
$ws.Range("L33").Value = "8b03b509-7093-40bc-8434-a8bb0e863b8b"
$ws.Range("L34").Value = "87692c93-012d-46df-9a34-c00c584ff1aa"
$ws.Range("L35").Value = "9c58e938-6860-48d2-811b-7a808fab0100"
$ws.Range("L36").Value = "35a25288-788d-4a44-8795-b32fe3cfa30e"
$ws.Range("L37").Value = "ca2d65c6-d4c3-4ce0-8feb-b509e67aab4a"
$ws.Range("L38").Value = "d52b5367-a67f-4c40-b3d2-3103b549415f"
$ws.Range("L39").Value = "c9806066-7233-45af-92e5-908edd1beacd"
$ws.Range("L40").Value = "98161028-ac3b-43f4-8786-5ca2aec3476f"
$ws.Range("L41").Value = "b3fc2e79-a20d-49a9-b5f7-65559ab2c2f1"
$ws.Range("L42").Value = "10917618-5003-4daa-a7fc-a28d72a7d49c"
$ws.Range("L43").Value = "888edf57-93cd-4561-b79e-66e20d7dabe0"
$ws.Range("L44").Value = "d1868a5a-c510-457b-b519-74860bf9847a"
$ws.Range("L45").Value = "bae8f082-014e-4ba1-aa1b-d487e83421a1"
$ws.Range("L46").Value = "18b3fc93-59e5-435e-b3c0-d42c7120c6bb"
$ws.Range("L47").Value = "d394acae-f653-4e2b-b653-5bc4df0dbc57"
$ws.Range("L48").Value = "c0e9df00-c3f3-49f3-be9e-108ad0df3235"
$ws.Range("L49").Value = "80de24d2-5d78-4adf-ad6b-298bac6c0b48"
$ws.Range("L50").Value = "0c6099a0-da26-4dd9-bb49-263672f992f5"
$ws.Range("L51").Value = "bb0d442b-2a9b-428f-a685-1fe48460f1a3"
$ws.Range("L52").Value = "fca78ffc-18eb-4374-a346-39061a1e26e2"
$ws.Range("L53").Value = "2217ed02-aee9-4123-bf65-54e608e2aaba"
$ws.Range("L54").Value = "9d1367c5-58df-40c8-8905-87f663ec6241"
$ws.Range("L55").Value = "eeacd526-85bd-4cdb-aa5d-a0142b9e7f49"
$ws.Range("L56").Value = "711f21eb-972b-43de-9f17-cb1358a6fe44"
$ws.Range("L57").Value = "037a23ed-b1be-43d8-b546-8cc483b876b5"
$ws.Range("L58").Value = "c23e2284-347a-45a6-bcc0-0a92c9c8a88c"
$ws.Range("L59").Value = "843d6518-76e5-40b8-a7b3-75ae9d5910c3"
$ws.Range("L60").Value = "89cf9c31-4a08-41ba-a938-a058de10fa9a"
$ws.Range("L61").Value = "85eb3468-c095-4f7d-bd07-b8effc1ce3f3"
$ws.Range("L62").Value = "ae44a254-c516-4736-b38b-49bf97a8de17"
$ws.Range("L63").Value = "d0d119bb-d7ad-47c2-9d75-5c4939270e5d"
$ws.Range("L64").Value = "fb4b8e18-fe4d-43c4-99cb-71c8df90c721"
$ws.Range("L65").Value = "7fc615f1-df1e-41a4-9605-b1d8f2887e9d"
$ws.Range("L66").Value = "71545c6e-cdd3-4d87-851b-5f7c473a52eb"
$ws.Range("L67").Value = "348cc385-9a23-41ff-9ef2-b37e2deaf323"
$ws.Range("L68").Value = "eb997444-68b6-469a-b2c5-1a52a78981b6"
$ws.Range("L69").Value = "06f7ad6a-821b-47b5-a4fa-511feb95c069"
$ws.Range("L70").Value = "5e560bff-b904-4d3c-966c-4fc5e73ac6f7"
$ws.Range("L71").Value = "c7a83500-2b26-41f8-a474-bed1a2563bab"
$ws.Range("L72").Value = "e6c3af4d-6e42-4647-b232-2866c22e60cb"
$ws.Range("L73").Value = "aa4b48d5-fd5b-4807-930d-e40e35a86364"
$ws.Range("L74").Value = "3c91856b-c6fc-4359-ae44-6a5b51866f5c"
$ws.Range("L75").Value = "4d6b304e-a1ce-417a-8488-3e0e6b1ce102"
$ws.Range("L76").Value = "df247289-e4be-44fe-b6fc-4407b1bf2e99"
$ws.Range("L77").Value = "ab699f75-21cb-4f69-8892-49127010a39d"
$ws.Range("L78").Value = "74c97b9c-79ba-4f4f-b66e-537ec3ad2346"
$ws.Range("L79").Value = "d81b5b77-5af0-4654-9705-74ee24a46fd4"
$ws.Range("L80").Value = "cef5bd3b-fb7e-4851-bf39-c3027d536e66"
$ws.Range("L81").Value = "a7c9f6ed-a681-4d52-ba32-01e6d018893e"
$ws.Range("L82").Value = "12a9a0e2-3daa-45be-81ab-1593dde148ec"
$ws.Range("L83").Value = "8779aafa-e18b-4af1-8749-0cdcf83f9358"
$ws.Range("L84").Value = "dc71351b-6d8d-4031-862c-34249f686a83"
$ws.Range("L85").Value = "6438da93-f364-4881-8f60-2e0b5eefd5c8"
$ws.Range("L86").Value = "4e9d069e-f036-499f-9c53-e2eb1624a31a"
$ws.Range("L87").Value = "a5a50acd-c8bb-4384-9817-870013dbbc46"
$ws.Range("L88").Value = "483d5484-1f8e-4a1c-984f-63b809a9588d"
$ws.Range("L89").Value = "216bb252-4d64-41d9-b191-80d6e5334a4c"
$ws.Range("L90").Value = "4edc032c-f374-4b49-be12-4ab481923798"
$ws.Range("L91").Value = "f108702e-b61d-4d35-bc47-68d78ff93bfd"
$ws.Range("L92").Value = "4fdff024-7e14-48b0-9254-e2b04893aad0"
$ws.Range("L93").Value = "cbafd002-fe39-4a69-8094-ca8a4ef111b5"
$ws.Range("L94").Value = "0cb25ecf-3771-4f9c-98a8-5234b577dad4"
$ws.Range("L95").Value = "f243d14f-def2-4469-b206-3e461d75dded"
$ws.Range("L96").Value = "a43f158f-7a3e-4801-8852-e1e2505035ca"
$ws.Range("L97").Value = "9d9ecca1-7776-4a00-99ce-ba50d95dfab2"
$ws.Range("L98").Value = "2c7a77bb-3c13-42f3-b8f4-7c85401235b5"
$ws.Range("L99").Value = "c4fa7a04-0f4e-45e9-acc1-4bf7b60322fb"
$ws.Range("L100").Value = "d86f090b-e141-4840-83ab-bb706a7b29f6"
$ws.Range("L101").Value = "648d0c48-32ce-41f8-81ee-8b1e1e2f1252"
$ws.Range("L102").Value = "bd65676e-654c-406c-b241-49b0c50e4441"
$ws.Range("L103").Value = "da1c08b4-5f36-4bcb-9c0f-1daf577732bb"
$ws.Range("L119").Value = "b506f62c-2f97-4388-945d-6fe4bf4f8875"
